$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    $escaped = $Text.Replace('"', '""')
    $Range.Formula = '="' + $escaped + '"'
}

Set-TextCell $ws.Range("D2") "93.867.90"
Set-TextCell $ws.Range("E2") "  -4.23%  "

Set-TextCell $ws.Range("D3") "3.429.68"
Set-TextCell $ws.Range("E3") "  +2.10%  "

Set-TextCell $ws.Range("D4") "1.00"
Set-TextCell $ws.Range("E4") "  -0.02%  "

Set-TextCell $ws.Range("D5") "235.80"
Set-TextCell $ws.Range("E5") "  -6.95%  "

Set-TextCell $ws.Range("D6") "638.07"
Set-TextCell $ws.Range("E6") "  -3.64%  "

Set-TextCell $ws.Range("D7") "1.42"
Set-TextCell $ws.Range("E7") "  -1.59%  "

Set-TextCell $ws.Range("D8") "0.393"
Set-TextCell $ws.Range("E8") "  -8.60%  "

Set-TextCell $ws.Range("E9") "  +0.10%  "

Set-TextCell $ws.Range("D10") "0.956"
Set-TextCell $ws.Range("E10") "  -6.92%  "

Set-TextCell $ws.Range("D11") "3.426.77"
Set-TextCell $ws.Range("E11") "  +2.08%  "

Set-TextCell $ws.Range("B12") "Avalanche"
Set-TextCell $ws.Range("C12") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws.Range("D12") "41.97"
Set-TextCell $ws.Range("E12") "  +0.43%  "

Set-TextCell $ws.Range("B13") "TRON"
Set-TextCell $ws.Range("C13") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws.Range("D13") "0.196"
Set-TextCell $ws.Range("E13") "  -5.65%  "

Set-TextCell $ws.Range("D14") "6.09"
Set-TextCell $ws.Range("E14") "  -0.91%  "

Set-TextCell $ws.Range("D15") "93.490.00"
Set-TextCell $ws.Range("E15") "  -4.26%  "

Set-TextCell $ws.Range("D16") "4.072.73"
Set-TextCell $ws.Range("E16") "  +2.35%  "

Set-TextCell $ws.Range("D17") "0.0000251"
Set-TextCell $ws.Range("E17") "  -2.34%  "

Set-TextCell $ws.Range("D18") "8.33"
Set-TextCell $ws.Range("E18") "  -7.01%  "

Set-TextCell $ws.Range("D19") "3.420.27"
Set-TextCell $ws.Range("E19") "  +2.44%  "

Set-TextCell $ws.Range("D20") "17.57"
Set-TextCell $ws.Range("E20") "  -2.64%  "

Set-TextCell $ws.Range("D21") "11.27"
Set-TextCell $ws.Range("E21") "  +3.74%  "

Set-TextCell $ws.Range("D22") "0.488"
Set-TextCell $ws.Range("E22") "  -12.21%  "

Set-TextCell $ws.Range("D23") "495.90"
Set-TextCell $ws.Range("E23") "  -3.53%  "

Set-TextCell $ws.Range("D24") "3.19"
Set-TextCell $ws.Range("E24") "  -5.33%  "

Set-TextCell $ws.Range("D25") "0.0000188"
Set-TextCell $ws.Range("E25") "  -6.81%  "

Set-TextCell $ws.Range("D26") "6.40"
Set-TextCell $ws.Range("E26") "  -4.27%  "

Set-TextCell $ws.Range("D27") "90.51"
Set-TextCell $ws.Range("E27") "  -7.26%  "

Set-TextCell $ws.Range("D28") "11.96"
Set-TextCell $ws.Range("E28") "  -3.06%  "

Set-TextCell $ws.Range("D29") "3.604.10"
Set-TextCell $ws.Range("E29") "  +2.13%  "

Set-TextCell $ws.Range("D30") "11.67"
Set-TextCell $ws.Range("E30") "  -0.50%  "

Set-TextCell $ws.Range("D31") "0.999"
Set-TextCell $ws.Range("E31") "  -0.13%  "

Set-TextCell $ws.Range("D32") "2.71"
Set-TextCell $ws.Range("E32") "  +5.46%  "

Set-TextCell $ws.Range("D33") "0.135"
Set-TextCell $ws.Range("E33") "  -8.83%  "

Set-TextCell $ws.Range("D34") "0.179"
Set-TextCell $ws.Range("E34") "  -6.34%  "

Set-TextCell $ws.Range("E35") "  -0.31%  "

Set-TextCell $ws.Range("D36") "29.94"
Set-TextCell $ws.Range("E36") "  +3.70%  "

Set-TextCell $ws.Range("D37") "0.556"
Set-TextCell $ws.Range("E37") "  -3.04%  "

Set-TextCell $ws.Range("D38") "541.99"
Set-TextCell $ws.Range("E38") "  +5.10%  "

Set-TextCell $ws.Range("D39") "7.57"
Set-TextCell $ws.Range("E39") "  -5.09%  "

Set-TextCell $ws.Range("D40") "1.44"
Set-TextCell $ws.Range("E40") "  -4.69%  "

Set-TextCell $ws.Range("E41") "  -0.11%  "

Set-TextCell $ws.Range("D42") "0.928"
Set-TextCell $ws.Range("E42") "  +8.92%  "

Set-TextCell $ws.Range("D43") "0.150"
Set-TextCell $ws.Range("E43") "  -1.88%  "

Set-TextCell $ws.Range("D44") "24.05"
Set-TextCell $ws.Range("E44") "  -1.57%  "

Set-TextCell $ws.Range("D45") "1.68"
Set-TextCell $ws.Range("E45") "  -3.33%  "

Set-TextCell $ws.Range("B46") "Filecoin"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D46") "5.54"
Set-TextCell $ws.Range("E46") "  -3.56%  "

Set-TextCell $ws.Range("B47") "VeChain"
Set-TextCell $ws.Range("C47") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D47") "0.0410"
Set-TextCell $ws.Range("E47") "  -5.73%  "

Set-TextCell $ws.Range("D48") "2.13"
Set-TextCell $ws.Range("E48") "  +4.83%  "

Set-TextCell $ws.Range("B49") "MantraDAO"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextCell $ws.Range("D49") "3.45"
Set-TextCell $ws.Range("E49") "  -5.45%  "

Set-TextCell $ws.Range("B50") "OKB"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D50") "52.77"
Set-TextCell $ws.Range("E50") "  -4.61%  "

Set-TextCell $ws.Range("D51") "3.16"
Set-TextCell $ws.Range("E51") "  -0.71%  "
